$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each triple is (row, column, new value) taken from the updated case results
# (case with 380 kV slack voltage of 1.02 pu instead of 1.05 pu)
$data = @(
    @(2, 2, 1.02),
    @(2, 3, 1.098088674980769),
    @(2, 4, 1.105694528146041),
    @(2, 5, 1.088209533992204),
    @(2, 6, 1.107066652842347),
    @(2, 9, 1.030414779499211),
    @(2, 10, 1.102882264957339),
    @(2, 11, 1.108303568967015),
    @(2, 12, 1.090862811102073),
    @(2, 13, 1.109672285910724),
    @(2, 14, 1.104448484275943),
    @(3, 2, 1.02),
    @(3, 3, 1.102959391417342),
    @(3, 4, 1.110449069066241),
    @(3, 5, 1.092523156267567),
    @(3, 6, 1.11166203794188),
    @(3, 9, 1.030642097196377),
    @(3, 10, 1.107403508153479),
    @(3, 11, 1.112871936866275),
    @(3, 12, 1.09498781784806),
    @(3, 13, 1.114082129908095),
    @(3, 14, 1.108976148156015),
    @(4, 2, 1.02),
    @(4, 3, 1.106070568209831),
    @(4, 4, 1.113485440477799),
    @(4, 5, 1.095276914843258),
    @(4, 6, 1.114595253736348),
    @(4, 9, 1.030783449794507),
    @(4, 10, 1.110289561085009),
    @(4, 11, 1.115788050440086),
    @(4, 12, 1.097619672142038),
    @(4, 13, 1.116895462012829),
    @(4, 14, 1.111866299613742),
    @(5, 2, 1.02),
    @(5, 3, 1.107369179490191),
    @(5, 4, 1.114752682223799),
    @(5, 5, 1.096425961107796),
    @(5, 6, 1.115819078032995),
    @(5, 9, 1.030841525027894),
    @(5, 10, 1.111493750476121),
    @(5, 11, 1.11700477651153),
    @(5, 12, 1.098717495439312),
    @(5, 13, 1.118068919930825),
    @(5, 14, 1.113072199091913),
    @(6, 2, 1.02),
    @(6, 3, 1.107586685528958),
    @(6, 4, 1.114964925524294),
    @(6, 5, 1.096618394098168),
    @(6, 6, 1.11602402809267),
    @(6, 9, 1.030851197667332),
    @(6, 10, 1.111695414971775),
    @(6, 11, 1.117208540134481),
    @(6, 12, 1.098901329126323),
    @(6, 13, 1.118265414920576),
    @(6, 14, 1.113274149974281),
    @(7, 2, 1.02),
    @(7, 3, 1.1060879564892),
    @(7, 4, 1.113502409299086),
    @(7, 5, 1.095292301952369),
    @(7, 6, 1.114611642624012),
    @(7, 9, 1.030784231069518),
    @(7, 10, 1.110305686856916),
    @(7, 11, 1.115804344120403),
    @(7, 12, 1.097634374717882),
    @(7, 13, 1.116911177785427),
    @(7, 14, 1.111882448286095),
    @(8, 2, 1.02),
    @(8, 3, 1.099743369558077),
    @(8, 4, 1.107309878670592),
    @(8, 5, 1.089675298802168),
    @(8, 6, 1.108628245297256),
    @(8, 9, 1.03049280532548),
    @(8, 10, 1.104418631465595),
    @(8, 11, 1.109855952608703),
    @(8, 12, 1.092264797645044),
    @(8, 13, 1.11117113243636),
    @(8, 14, 1.105987032600865),
    @(9, 2, 1.02),
    @(9, 3, 1.088235870116299),
    @(9, 4, 1.096073583794844),
    @(9, 5, 1.079475141734748),
    @(9, 6, 1.097759621190553),
    @(9, 9, 1.029934201343458),
    @(9, 10, 1.093726086193293),
    @(9, 11, 1.099051890669819),
    @(9, 12, 1.082502128622125),
    @(9, 13, 1.100733072756442),
    @(9, 14, 1.095279302687823),
    @(10, 2, 1.02),
    @(10, 3, 1.080318567965686),
    @(10, 4, 1.088339930342042),
    @(10, 5, 1.072449045978129),
    @(10, 6, 1.090271188024423),
    @(10, 9, 1.02952982200379),
    @(10, 10, 1.086359419506987),
    @(10, 11, 1.091608387784106),
    @(10, 12, 1.075769278950562),
    @(10, 13, 1.093533465656913),
    @(10, 14, 1.087902174489854),
    @(11, 2, 1.02),
    @(11, 3, 1.076825984059971),
    @(11, 4, 1.084927698830576),
    @(11, 5, 1.069347659027814),
    @(11, 6, 1.086965325642246),
    @(11, 9, 1.029346738060513),
    @(11, 10, 1.083107354532027),
    @(11, 11, 1.08832241755796),
    @(11, 12, 1.072795385631259),
    @(11, 13, 1.090353235146621),
    @(11, 14, 1.084645491209613),
    @(12, 2, 1.02),
    @(12, 3, 1.075518468687145),
    @(12, 4, 1.083650168936316),
    @(12, 5, 1.068186305326908),
    @(12, 6, 1.085727348650377),
    @(12, 9, 1.029277496575626),
    @(12, 10, 1.08188952420399),
    @(12, 11, 1.087091894033819),
    @(12, 12, 1.071681478781001),
    @(12, 13, 1.089162021821121),
    @(12, 14, 1.083425931422822),
    @(13, 2, 1.02),
    @(13, 3, 1.075799406535542),
    @(13, 4, 1.08392466827037),
    @(13, 5, 1.068435851405535),
    @(13, 6, 1.085993361578363),
    @(13, 9, 1.029292405645373),
    @(13, 10, 1.082151208138834),
    @(13, 11, 1.087356305260222),
    @(13, 12, 1.071920843066231),
    @(13, 13, 1.089417999114164),
    @(13, 14, 1.08368798697887),
    @(14, 2, 1.02),
    @(14, 3, 1.076718116077813),
    @(14, 4, 1.084822306536692),
    @(14, 5, 1.06925185505448),
    @(14, 6, 1.08686320184591),
    @(14, 9, 1.029341039968344),
    @(14, 10, 1.083006892704245),
    @(14, 11, 1.088220908555298),
    @(14, 12, 1.072703501717904),
    @(14, 13, 1.090254974764436),
    @(14, 14, 1.084544886714512),
    @(15, 2, 1.02),
    @(15, 3, 1.077282792945092),
    @(15, 4, 1.085374019554267),
    @(15, 5, 1.06975336633653),
    @(15, 6, 1.087397793703342),
    @(15, 9, 1.029370840281826),
    @(15, 10, 1.083532784575155),
    @(15, 11, 1.088752282273774),
    @(15, 12, 1.07318448036326),
    @(15, 13, 1.090769331086599),
    @(15, 14, 1.0850715254122),
    @(16, 2, 1.02),
    @(16, 3, 1.080548957641251),
    @(16, 4, 1.088565005994887),
    @(16, 5, 1.072653589680603),
    @(16, 6, 1.090489209474298),
    @(16, 9, 1.029541801397023),
    @(16, 10, 1.086573893424341),
    @(16, 11, 1.091825098179554),
    @(16, 12, 1.075965373245882),
    @(16, 13, 1.093743162400807),
    @(16, 14, 1.088116952984771),
    @(17, 2, 1.02),
    @(17, 3, 1.082580133877988),
    @(17, 4, 1.090549257168752),
    @(17, 5, 1.074456678827848),
    @(17, 6, 1.092411062709513),
    @(17, 9, 1.029646878008233),
    @(17, 10, 1.088464479390592),
    @(17, 11, 1.093735400523138),
    @(17, 12, 1.077693756337872),
    @(17, 13, 1.095591419120876),
    @(17, 14, 1.090010223799948),
    @(18, 2, 1.02),
    @(18, 3, 1.083758715370909),
    @(18, 4, 1.091700547276454),
    @(18, 5, 1.075502728163057),
    @(18, 6, 1.093525974183618),
    @(18, 9, 1.029707399540612),
    @(18, 10, 1.089561256398316),
    @(18, 11, 1.09484361680003),
    @(18, 12, 1.078696279815682),
    @(18, 13, 1.096663456907097),
    @(18, 14, 1.091108558356825),
    @(19, 2, 1.02),
    @(19, 3, 1.084159550915401),
    @(19, 4, 1.092092090184106),
    @(19, 5, 1.075858457998562),
    @(19, 6, 1.093905115377405),
    @(19, 9, 1.029727906641241),
    @(19, 10, 1.089934231626868),
    @(19, 11, 1.095220482299439),
    @(19, 12, 1.079037176388342),
    @(19, 13, 1.097027988042648),
    @(19, 14, 1.091482063252982),
    @(20, 2, 1.02),
    @(20, 3, 1.082362849863556),
    @(20, 4, 1.09033699947703),
    @(20, 5, 1.074263813533468),
    @(20, 6, 1.092205497862282),
    @(20, 9, 1.02963568396303),
    @(20, 10, 1.088262258530132),
    @(20, 11, 1.093531070615051),
    @(20, 12, 1.077508901148021),
    @(20, 13, 1.095393744814845),
    @(20, 14, 1.08980771576267),
    @(21, 2, 1.02),
    @(21, 3, 1.07644786553783),
    @(21, 4, 1.084558257006801),
    @(21, 5, 1.069011824793765),
    @(21, 6, 1.086607336834616),
    @(21, 9, 1.029326752793176),
    @(21, 10, 1.082755191645702),
    @(21, 11, 1.087966583926296),
    @(21, 12, 1.072473288138369),
    @(21, 13, 1.090008784739351),
    @(21, 14, 1.084292828211595),
    @(22, 2, 1.02),
    @(22, 3, 1.072669489033344),
    @(22, 4, 1.080866353138239),
    @(22, 5, 1.065655270269163),
    @(22, 6, 1.083029226420204),
    @(22, 9, 1.029125345077349),
    @(22, 10, 1.079235305608551),
    @(22, 11, 1.08441001929142),
    @(22, 12, 1.069253305272066),
    @(22, 13, 1.086565298503553),
    @(22, 14, 1.080767943532534),
    @(23, 2, 1.02),
    @(23, 3, 1.074678298860124),
    @(23, 4, 1.08282924027204),
    @(23, 5, 1.067439973004005),
    @(23, 6, 1.084931760255443),
    @(23, 9, 1.029232807480627),
    @(23, 10, 1.081106882560087),
    @(23, 11, 1.086301096274631),
    @(23, 12, 1.070965554413382),
    @(23, 13, 1.088396406348006),
    @(23, 14, 1.082642178338015),
    @(24, 2, 1.02),
    @(24, 3, 1.082461050182496),
    @(24, 4, 1.090432928373324),
    @(24, 5, 1.074350978520649),
    @(24, 6, 1.092298402498641),
    @(24, 9, 1.029640744440544),
    @(24, 10, 1.088353651844062),
    @(24, 11, 1.093623417104933),
    @(24, 12, 1.077592446563426),
    @(24, 13, 1.09548308388956),
    @(24, 14, 1.089899238865589),
    @(25, 2, 1.02),
    @(25, 3, 1.091252183279084),
    @(25, 4, 1.099019332539834),
    @(25, 5, 1.082150217529299),
    @(25, 6, 1.100610351290175),
    @(25, 9, 1.030084121830008),
    @(25, 10, 1.096530533979314),
    @(25, 11, 1.101885593418699),
    @(25, 12, 1.085063868126499),
    @(25, 13, 1.103472223830235),
    @(25, 14, 1.098087733111375)
)

foreach ($item in $data) {
    $ws.Cells.Item($item[0], $item[1]).Value = $item[2]
}
